$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values for rows that were reshuffled (matches target diff row-by-row)
$ws.Range("A2").Value = 111896639
$ws.Range("Q2").Value = 575089.384229039
$ws.Range("R2").Value = 6703379.745088123
$ws.Range("A3").Value = 111884471
$ws.Range("B3").Value = 88899
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 3286
$ws.Range("F3").Value = 'Flattoppad klubbsvamp'
$ws.Range("G3").Value = 'Clavariadelphus truncatus'
$ws.Range("H3").Value = '(Quél.) Donk'
$ws.Range("P3").Value = 'Kalkberget (Kalkberget), Gstr'
$ws.Range("Q3").Value = 575020.8210917887
$ws.Range("R3").Value = 6703397.074168184
$ws.Range("AW3").Value = 'Patric Engfeldt'
$ws.Range("AX3").Value = 'Patric Engfeldt'
$ws.Range("A4").Value = 111896690
$ws.Range("B4").Value = 90687
$ws.Range("E4").Value = 5964
$ws.Range("F4").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G4").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H4").Value = '(L.:Fr.) P.Karst.'
$ws.Range("Q4").Value = 575060.2881161601
$ws.Range("R4").Value = 6703376.67477417
$ws.Range("A5").Value = 111896634
$ws.Range("B5").Value = 90332
$ws.Range("E5").Value = 4769
$ws.Range("F5").Value = 'Svavelriska'
$ws.Range("G5").Value = 'Lactarius scrobiculatus'
$ws.Range("H5").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P5").Value = 'Kratte masugn, Gstr'
$ws.Range("Q5").Value = 575048.3395925189
$ws.Range("R5").Value = 6703452.413791304
$ws.Range("AW5").Value = 'Philipp Weiss'
$ws.Range("AX5").Value = 'Philipp Weiss'
$ws.Range("A6").Value = 111896653
$ws.Range("B6").Value = 89183
$ws.Range("E6").Value = 3215
$ws.Range("F6").Value = 'Rödgul trumpetsvamp'
$ws.Range("G6").Value = 'Craterellus lutescens'
$ws.Range("H6").Value = '(Fr.) Fr.'
$ws.Range("P6").Value = 'Kratte masugn, Gstr'
$ws.Range("Q6").Value = 575075.050630242
$ws.Range("R6").Value = 6703403.625642136
$ws.Range("AW6").Value = 'Philipp Weiss'
$ws.Range("AX6").Value = 'Philipp Weiss'
$ws.Range("A7").Value = 111896641
$ws.Range("B7").Value = 90332
$ws.Range("E7").Value = 4769
$ws.Range("F7").Value = 'Svavelriska'
$ws.Range("G7").Value = 'Lactarius scrobiculatus'
$ws.Range("H7").Value = '(Scop.:Fr.) Fr.'
$ws.Range("Q7").Value = 575021.3626164712
$ws.Range("R7").Value = 6703370.933926445
$ws.Range("A8").Value = 111896633
$ws.Range("B8").Value = 90332
$ws.Range("E8").Value = 4769
$ws.Range("F8").Value = 'Svavelriska'
$ws.Range("G8").Value = 'Lactarius scrobiculatus'
$ws.Range("H8").Value = '(Scop.:Fr.) Fr.'
$ws.Range("Q8").Value = 575100.4050603262
$ws.Range("R8").Value = 6703444.118284944
$ws.Range("A9").Value = 111896644
$ws.Range("Q9").Value = 575036.4083237475
$ws.Range("R9").Value = 6703431.936489306
$ws.Range("A11").Value = 111896640
$ws.Range("Q11").Value = 575025.3556637274
$ws.Range("R11").Value = 6703369.042946251
$ws.Range("A12").Value = 111896643
$ws.Range("B12").Value = 90332
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 4769
$ws.Range("F12").Value = 'Svavelriska'
$ws.Range("G12").Value = 'Lactarius scrobiculatus'
$ws.Range("H12").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P12").Value = 'Kratte masugn, Gstr'
$ws.Range("Q12").Value = 575038.7114136803
$ws.Range("R12").Value = 6703416.194821274
$ws.Range("AW12").Value = 'Philipp Weiss'
$ws.Range("AX12").Value = 'Philipp Weiss'
$ws.Range("A13").Value = 111896638
$ws.Range("Q13").Value = 575087.1320314853
$ws.Range("R13").Value = 6703393.020834555
$ws.Range("A14").Value = 111896637
$ws.Range("B14").Value = 90332
$ws.Range("E14").Value = 4769
$ws.Range("F14").Value = 'Svavelriska'
$ws.Range("G14").Value = 'Lactarius scrobiculatus'
$ws.Range("H14").Value = '(Scop.:Fr.) Fr.'
$ws.Range("Q14").Value = 575088.0587098968
$ws.Range("R14").Value = 6703396.00058554
$ws.Range("A15").Value = 111896654
$ws.Range("B15").Value = 89183
$ws.Range("E15").Value = 3215
$ws.Range("F15").Value = 'Rödgul trumpetsvamp'
$ws.Range("G15").Value = 'Craterellus lutescens'
$ws.Range("H15").Value = '(Fr.) Fr.'
$ws.Range("Q15").Value = 575072.6962527435
$ws.Range("R15").Value = 6703421.833381963
$ws.Range("A16").Value = 111896655
$ws.Range("B16").Value = 89183
$ws.Range("E16").Value = 3215
$ws.Range("F16").Value = 'Rödgul trumpetsvamp'
$ws.Range("G16").Value = 'Craterellus lutescens'
$ws.Range("H16").Value = '(Fr.) Fr.'
$ws.Range("Q16").Value = 575104.6742508161
$ws.Range("R16").Value = 6703428.910891063
$ws.Range("A17").Value = 111896636
$ws.Range("Q17").Value = 575108.85141061
$ws.Range("R17").Value = 6703418.142308297
$ws.Range("A18").Value = 111883983
$ws.Range("P18").Value = 'Kalkberget (Kalkberget), Gstr'
$ws.Range("Q18").Value = 575058.3527020445
$ws.Range("R18").Value = 6703446.206921679
$ws.Range("AW18").Value = 'Patric Engfeldt'
$ws.Range("AX18").Value = 'Patric Engfeldt'
$ws.Range("A19").Value = 111896635
$ws.Range("Q19").Value = 575037.2974304935
$ws.Range("R19").Value = 6703389.027347369
$ws.Range("A20").Value = 111884133
$ws.Range("B20").Value = 88899
$ws.Range("D20").Value = 'NT'
$ws.Range("E20").Value = 3286
$ws.Range("F20").Value = 'Flattoppad klubbsvamp'
$ws.Range("G20").Value = 'Clavariadelphus truncatus'
$ws.Range("H20").Value = '(Quél.) Donk'
$ws.Range("P20").Value = 'Kalkberget (Kalkberget), Gstr'
$ws.Range("Q20").Value = 575059.034285416
$ws.Range("R20").Value = 6703389.477814267
$ws.Range("AW20").Value = 'Patric Engfeldt'
$ws.Range("AX20").Value = 'Patric Engfeldt'
$ws.Range("A22").Value = 111884093
$ws.Range("B22").Value = 98535
$ws.Range("D22").Value = 'LC'
$ws.Range("E22").Value = 222498
$ws.Range("F22").Value = 'Blåsippa'
$ws.Range("G22").Value = 'Hepatica nobilis'
$ws.Range("H22").Value = 'Schreb.'
$ws.Range("P22").Value = 'Kopparåsen (Kopparåsen), Gstr'
$ws.Range("Q22").Value = 575065.9914513066
$ws.Range("R22").Value = 6703387.648325931

# Re-create empty placeholder cells that moved into these rows
$ws.Range("K3").NumberFormat = "General"
$ws.Range("AF4").NumberFormat = "General"
$ws.Range("K18").NumberFormat = "General"
$ws.Range("K20").NumberFormat = "General"

# Remove empty placeholder cells that moved out of these rows
$ws.Range("K5").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("AF14").ClearContents()
